$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1144.7778
$ws.Range("I2").Value = 1144.7778
$ws.Range("K2").Value = 1144.7778
$ws.Range("M2").Value = -1031.7778
$ws.Range("H18").Value = 2422.1177
$ws.Range("I18").Value = 2345.6
$ws.Range("J18").Value = 2996
$ws.Range("K18").Value = 2345.6
$ws.Range("L18").Value = 2996
$ws.Range("M18").Value = -2061.6
$ws.Range("N18").Value = -3564
$ws.Range("H43").Value = 2125.9285
$ws.Range("J43").Value = 2234
$ws.Range("L43").Value = 2234
$ws.Range("N43").Value = -2372
$ws.Range("H94").Value = 50128996
$ws.Range("I94").Value = 71437140
$ws.Range("J94").Value = 410002
$ws.Range("K94").Value = 71437140
$ws.Range("L94").Value = 410002
$ws.Range("M94").Value = -71436689
$ws.Range("N94").Value = -410904
$ws.Range("H106").Value = 5547
$ws.Range("I106").Value = 4559.8667
$ws.Range("K106").Value = 4559.8667
$ws.Range("M106").Value = -3928.8667
$ws.Range("H138").Value = 9895.821
$ws.Range("J138").Value = 9597.388999999999
$ws.Range("L138").Value = 28792.167
$ws.Range("N138").Value = -39072.167
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 5604.607
$ws.Range("I141").Value = 4406.5454
$ws.Range("J141").Value = 9997.5
$ws.Range("K141").Value = 13219.6362
$ws.Range("L141").Value = 29992.5
$ws.Range("M141").Value = -8039.636200000001
$ws.Range("N141").Value = -40352.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 29498280
$ws.Range("I74").Value = 38573136
$ws.Range("K74").Value = 38573136
$ws.Range("M74").Value = -38572262
$ws.Range("H77").Value = 29498280
$ws.Range("I77").Value = 38573136
$ws.Range("K77").Value = 192865680
$ws.Range("M77").Value = -192861312
$ws.Range("H82").Value = 55000
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55722
$ws.Range("H85").Value = 55000
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57496
$ws.Range("H102").Value = 5490.18
$ws.Range("I102").Value = 3701.5642
$ws.Range("K102").Value = 3701.5642
$ws.Range("M102").Value = -2079.5642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4724.8184
$ws.Range("I20").Value = 2216.3333
$ws.Range("K20").Value = 2216.3333
$ws.Range("M20").Value = -1969.3333
$ws.Range("H86").Value = 5623.7915
$ws.Range("I86").Value = 6566.5
$ws.Range("K86").Value = 6566.5
$ws.Range("M86").Value = -5443.5
$ws.Range("H89").Value = 5623.7915
$ws.Range("I89").Value = 6566.5
$ws.Range("K89").Value = 32832.5
$ws.Range("M89").Value = -27216.5
$ws.Range("H94").Value = 12264.857
$ws.Range("I94").Value = 17539
$ws.Range("K94").Value = 17539
$ws.Range("M94").Value = -17088
$ws.Range("H95").Value = 31524.8
$ws.Range("J95").Value = 31524.8
$ws.Range("L95").Value = 31524.8
$ws.Range("N95").Value = -37016.8
$ws.Range("H97").Value = 11319.5
$ws.Range("I97").Value = 11319.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 11319.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -10328.5
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 12868
$ws.Range("I99").Value = 14486.571
$ws.Range("K99").Value = 14486.571
$ws.Range("M99").Value = -12988.571
$ws.Range("H106").Value = 49700
$ws.Range("J106").Value = 49700
$ws.Range("L106").Value = 49700
$ws.Range("N106").Value = -52224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5560.024
$ws.Range("I58").Value = 6764.2915
$ws.Range("J58").Value = 3954.3333
$ws.Range("K58").Value = 6764.2915
$ws.Range("L58").Value = 3954.3333
$ws.Range("M58").Value = -6561.2915
$ws.Range("N58").Value = -4360.3333
$ws.Range("H105").Value = 151741.42
$ws.Range("I105").Value = 192126.36
$ws.Range("J105").Value = 3663.3333
$ws.Range("K105").Value = 192126.36
$ws.Range("L105").Value = 3663.3333
$ws.Range("M105").Value = -190379.36
$ws.Range("N105").Value = -7157.3333
$ws.Range("H107").Value = 16676.715
$ws.Range("I107").Value = 22277.4
$ws.Range("J107").Value = 2675
$ws.Range("K107").Value = 22277.4
$ws.Range("L107").Value = 2675
$ws.Range("M107").Value = -20357.4
$ws.Range("N107").Value = -6515
$ws.Range("H132").Value = 33386666
$ws.Range("I132").Value = 55576110
$ws.Range("J132").Value = 102500
$ws.Range("K132").Value = 166728330
$ws.Range("L132").Value = 307500
$ws.Range("M132").Value = -166725800
$ws.Range("N132").Value = -312560
$ws.Range("H136").Value = 5560.024
$ws.Range("I136").Value = 6764.2915
$ws.Range("J136").Value = 3954.3333
$ws.Range("K136").Value = 20292.8745
$ws.Range("L136").Value = 11862.9999
$ws.Range("M136").Value = -17742.8745
$ws.Range("N136").Value = -16962.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50123570
$ws.Range("I4").Value = 73506664
$ws.Range("K4").Value = 220519992
$ws.Range("M4").Value = -220519880
$ws.Range("H107").Value = 1930.3889
$ws.Range("J107").Value = 2059.25
$ws.Range("L107").Value = 6177.75
$ws.Range("N107").Value = -10017.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5232.788
$ws.Range("I70").Value = 5242.0713
$ws.Range("K70").Value = 5242.0713
$ws.Range("M70").Value = -4972.0713
$ws.Range("H73").Value = 5232.788
$ws.Range("I73").Value = 5242.0713
$ws.Range("K73").Value = 5242.0713
$ws.Range("M73").Value = -4306.0713
$ws.Range("H126").Value = 9907.973
$ws.Range("I126").Value = 11872.5
$ws.Range("J126").Value = 7943.4443
$ws.Range("K126").Value = 35617.5
$ws.Range("L126").Value = 23830.3329
$ws.Range("M126").Value = -33147.5
$ws.Range("N126").Value = -28770.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50676
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52340
$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52773

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H45").Value = 17313
$ws.Range("J45").Value = 17313
$ws.Range("L45").Value = 17313
$ws.Range("N45").Value = -18295
$ws.Range("H81").Value = 11882.632
$ws.Range("I81").Value = 12181.667
$ws.Range("J81").Value = 6500
$ws.Range("K81").Value = 24363.334
$ws.Range("L81").Value = 13000
$ws.Range("M81").Value = -23302.334
$ws.Range("N81").Value = -15122
$ws.Range("H84").Value = 11882.632
$ws.Range("I84").Value = 12181.667
$ws.Range("J84").Value = 6500
$ws.Range("K84").Value = 121816.67
$ws.Range("L84").Value = 65000
$ws.Range("M84").Value = -116512.67
$ws.Range("N84").Value = -75608
$ws.Range("H96").Value = 27084972
$ws.Range("I96").Value = 12501561
$ws.Range("J96").Value = 41668380
$ws.Range("K96").Value = 12501561
$ws.Range("L96").Value = 41668380
$ws.Range("M96").Value = -12500188
$ws.Range("N96").Value = -41671126
$ws.Range("H100").Value = 38860.645
$ws.Range("I100").Value = 22049.908
$ws.Range("K100").Value = 44099.816
$ws.Range("M100").Value = -43558.816
$ws.Range("H126").Value = 21989.652
$ws.Range("I126").Value = 38818
$ws.Range("K126").Value = 116454
$ws.Range("M126").Value = -113984
$ws.Range("H130").Value = 64999.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 64999.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 64999.5
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -75039.5
